$d = $word.ActiveDocument
$t = $d.Tables(1)

function Set-CellXml($table, $row, $col, $newLines) {
    $c = $table.Cell($row, $col)
    $r = $c.Range

    $parts = @()
    for ($i = 0; $i -lt $newLines.Length; $i++) {
        $line = $newLines[$i]
        $escaped = $line -replace '&','&amp;' -replace '<','&lt;' -replace '>','&gt;'
        if ($line -ne $line.Trim()) {
            $parts += "<w:t xml:space=`"preserve`">$escaped</w:t>"
        } else {
            $parts += "<w:t>$escaped</w:t>"
        }
    }
    $runInner = [string]::Join("<w:br/>", $parts)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:sz w:val="32"/></w:rPr>' + $runInner + '</w:r></w:p>'
    $r.InsertXML($xml)
}

Set-CellXml $t 1 1 @("71 x 59", "  5    9", "  ----", "7|    |", "1|    |")
Set-CellXml $t 1 2 @("45 x 21", "  2    1", "  ----", "4|    |", "5|    |")
Set-CellXml $t 1 3 @("26 x 47", "  4    7", "  ----", "2|    |", "6|    |")
Set-CellXml $t 2 1 @("86 x 35", "  3    5", "  ----", "8|    |", "6|    |")
Set-CellXml $t 2 2 @("24 x 15", "  1    5", "  ----", "2|    |", "4|    |")
Set-CellXml $t 2 3 @("61 x 36", "  3    6", "  ----", "6|    |", "1|    |")
Set-CellXml $t 3 1 @("65 x 57", "  5    7", "  ----", "6|    |", "5|    |")
Set-CellXml $t 3 2 @("72 x 71", "  7    1", "  ----", "7|    |", "2|    |")
Set-CellXml $t 3 3 @("90 x 75", "  7    5", "  ----", "9|    |", "0|    |")
Set-CellXml $t 4 1 @("44 x 31", "  3    1", "  ----", "4|    |", "4|    |")
Set-CellXml $t 4 2 @("45 x 88", "  8    8", "  ----", "4|    |", "5|    |")
Set-CellXml $t 4 3 @("44 x 41", "  4    1", "  ----", "4|    |", "4|    |")
Set-CellXml $t 5 1 @("87 x 76", "  7    6", "  ----", "8|    |", "7|    |")
Set-CellXml $t 5 2 @("68 x 50", "  5    0", "  ----", "6|    |", "8|    |")
Set-CellXml $t 5 3 @("80 x 10", "  1    0", "  ----", "8|    |", "0|    |")

Write-Host "Done"
